$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: semantic type URIs
$ws.Range("D2").Value2 = "iaest-measure:residencia-comarca-nombre"
$ws.Range("E2").Value2 = "iaest-measure:sexo"
$ws.Range("F2").Value2 = "iaest-measure:residencia-provincia-nombre"
$ws.Range("H2").Value2 = "iaest-measure:edad-grupos-quinquenales"
$ws.Range("I2").Value2 = "iaest-measure:residencia-ccaa-nombre"

# Row 3: dim -> medida
$ws.Range("D3").Value2 = "medida"
$ws.Range("E3").Value2 = "medida"
$ws.Range("F3").Value2 = "medida"
$ws.Range("H3").Value2 = "medida"
$ws.Range("I3").Value2 = "medida"

# Row 4: datatype -> xsd:int
$ws.Range("D4").Value2 = "xsd:int"
$ws.Range("E4").Value2 = "xsd:int"
$ws.Range("F4").Value2 = "xsd:int"
$ws.Range("H4").Value2 = "xsd:int"
$ws.Range("I4").Value2 = "xsd:int"

# Row 5: remove mapping file references that are no longer needed
$ws.Range("E5").Clear()
$ws.Range("H5").Clear()
$ws.Range("I5").Clear()
